$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "27.993.17"
$ws.Range('E2').Value = "  +3.45%  "
$ws.Range('D3').Value = "1.727.43"
$ws.Range('E3').Value = "  +3.13%  "
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = "  -0.18%  "
$ws.Range('D5').Value = "'219.24"
$ws.Range('E5').Value = "  +1.96%  "
$ws.Range('E6').Value = "  +1.46%  "
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = "  -0.13%  "
$ws.Range('D8').Value = "'24.25"
$ws.Range('E8').Value = "  +14.40%  "
$ws.Range('E9').Value = "  +3.78%  "
$ws.Range('E10').Value = "  +2.37%  "
$ws.Range('D11').Value = "'0.0900"
$ws.Range('E11').Value = "  +1.97%  "
$ws.Range('D12').Value = "1.971.26"
$ws.Range('E12').Value = "  +3.18%  "
$ws.Range('D13').Value = "1.719.48"
$ws.Range('E13').Value = "  +2.91%  "
$ws.Range('E14').Value = "  +3.72%  "
$ws.Range('D15').Value = "'0.567"
$ws.Range('E15').Value = "  +6.32%  "
$ws.Range('D16').Value = "'67.95"
$ws.Range('E16').Value = "  +2.97%  "
$ws.Range('D17').Value = "27.936.67"
$ws.Range('E17').Value = "  +3.34%  "
$ws.Range('D18').Value = "'243.92"
$ws.Range('D19').Value = "0.0₃0758"
$ws.Range('E19').Value = "  +2.54%  "
$ws.Range('D20').Value = "'7.92"
$ws.Range('E20').Value = "  -2.86%  "
$ws.Range('E21').Value = "  -0.11%  "
$ws.Range('E22').Value = "  +4.63%  "
$ws.Range('D23').Value = "'9.79"
$ws.Range('E23').Value = "  +4.86%  "
$ws.Range('E24').Value = "  +0.31%  "
$ws.Range('D25').Value = "'149.59"
$ws.Range('E25').Value = "  +2.33%  "
$ws.Range('D26').Value = "'7.54"
$ws.Range('E26').Value = "  +4.43%  "
$ws.Range('D27').Value = "'16.84"
$ws.Range('E27').Value = "  +3.00%  "
$ws.Range('E28').Value = "  +2.04%  "
$ws.Range('E30').Value = "  +2.74%  "
$ws.Range('D31').Value = "'1.20"
$ws.Range('E31').Value = "  +2.20%  "
$ws.Range('E32').Value = "  +3.11%  "
$ws.Range('D33').Value = "'3.28"
$ws.Range('E33').Value = "  +3.25%  "
$ws.Range('D34').Value = "1.488.19"
$ws.Range('E34').Value = "  -3.45%  "
$ws.Range('E35').Value = "  -1.67%  "
$ws.Range('D36').Value = "'0.962"
$ws.Range('E36').Value = "  +4.49%  "
$ws.Range('E37').Value = "  +2.39%  "
$ws.Range('E38').Value = "  +0.61%  "
$ws.Range('D39').Value = "'0.0176"
$ws.Range('E39').Value = "  +0.98%  "
$ws.Range('E40').Value = "  +1.02%  "
$ws.Range('D41').Value = "'71.52"
$ws.Range('E41').Value = "  +5.66%  "
$ws.Range('D43').Value = "'0.999"
$ws.Range('E43').Value = "  -0.11%  "
$ws.Range('D44').Value = "'2.30"
$ws.Range('E44').Value = "  +1.66%  "
$ws.Range('D45').Value = "1.875.23"
$ws.Range('E45').Value = "  +3.24%  "
$ws.Range('D46').Value = "'0.793"
$ws.Range('E46').Value = "  +1.05%  "
$ws.Range('D47').Value = "'1.78"
$ws.Range('E47').Value = "  +13.83%  "
$ws.Range('D48').Value = "'91.81"
$ws.Range('E48').Value = "  +1.07%  "
$ws.Range('E49').Value = "  +4.53%  "
$ws.Range('E50').Value = "  +1.47%  "
$ws.Range('D51').Value = "'8.24"
$ws.Range('E51').Value = "  +2.59%  "
